$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$data = @(
    @("SingleUseId297", "displayMeas",  "Left",   "LTR", "Graph with the measurement values from stamps mode"),
    @("SingleUseId298", "Default",      "Center", "LTR", "Frequency measurement 1"),
    @("SingleUseId299", "Default",      "Center", "LTR", "Frequency measurement 2"),
    @("SingleUseId300", "Default",      "Center", "LTR", "Frequency measurement 3"),
    @("SingleUseId301", "Default",      "Center", "LTR", "Frequency measurement 4"),
    @("SingleUseId302", "Default",      "Center", "LTR", "Frequency measurement 5"),
    @("SingleUseId303", "Default",      "Center", "LTR", "Frequency measurement 6"),
    @("SingleUseId304", "Default",      "Center", "LTR", "Frequency measurement 7"),
    @("SingleUseId305", "Default",      "Center", "LTR", "Frequency measurement 8"),
    @("SingleUseId306", "displayMeas",  "Left",   "LTR", "Graph - Frequency Measurement <value>"),
    @("SingleUseId307", "displayLabel", "Left",   "LTR", "<>"),
    @("SingleUseId308", "displayLabel", "Right",  "LTR", "<>")
)

$startRow = 237
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
    $ws.Cells.Item($row, 6).Value = $data[$i][4]
}
